# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# cells for the f3dc8a13-... row on the zh-cn and de-de report sheets,
# regenerating a fresh handback report timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 16:35:40"
$wsZhCn.Range("H4").Value = "2016-03-12 16:36:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 16:35:43"
$wsDeDe.Range("H4").Value = "2016-03-12 16:36:12"
